$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "28.009.55"
$ws.Range("E2").Value2 = "  -1.91%  "
$ws.Range("D3").Value2 = "1.884.90"
$ws.Range("E3").Value2 = "  -1.44%  "
$ws.Range("D4").Formula = "'1.003"
$ws.Range("E4").Value2 = "  +0.19%  "
$ws.Range("D5").Formula = "'313.08"
$ws.Range("E5").Value2 = "  -0.67%  "
$ws.Range("D6").Formula = "'1.002"
$ws.Range("E6").Value2 = "  +0.09%  "
$ws.Range("D7").Formula = "'0.4995"
$ws.Range("E7").Value2 = "  -4.45%  "
$ws.Range("D8").Formula = "'0.3869"
$ws.Range("E8").Value2 = "  -2.29%  "
$ws.Range("D9").Formula = "'0.09132"
$ws.Range("E9").Value2 = "  -5.89%  "
$ws.Range("D10").Formula = "'1.123"
$ws.Range("E10").Value2 = "  -2.59%  "
$ws.Range("D11").Formula = "'41.68"
$ws.Range("E11").Value2 = "  -0.67%  "
$ws.Range("D12").Formula = "'6.324"
$ws.Range("E12").Value2 = "  -3.27%  "
$ws.Range("D13").Formula = "'20.70"
$ws.Range("E13").Value2 = "  -2.30%  "
$ws.Range("D14").Value2 = "1.884.22"
$ws.Range("D15").Formula = "'7.280"
$ws.Range("E15").Value2 = "  -3.64%  "
$ws.Range("D16").Formula = "'1.003"
$ws.Range("E16").Value2 = "  +0.20%  "
$ws.Range("D17").Formula = "'0.00001103"
$ws.Range("E17").Value2 = "  -3.04%  "
$ws.Range("D18").Formula = "'91.53"
$ws.Range("E18").Value2 = "  -3.22%  "
$ws.Range("D19").Formula = "'0.06625"
$ws.Range("E19").Value2 = "  -0.38%  "
$ws.Range("D20").Formula = "'17.81"
$ws.Range("E20").Value2 = "  -2.15%  "
$ws.Range("E21").Value2 = "  +0.05%  "
$ws.Range("D22").Formula = "'6.204"
$ws.Range("E22").Value2 = "  -2.05%  "
$ws.Range("D23").Value2 = "28.080.46"
$ws.Range("E23").Value2 = "  -1.97%  "
$ws.Range("D24").Formula = "'11.32"
$ws.Range("E24").Value2 = "  -1.61%  "
$ws.Range("E25").Value2 = "  +0.48%  "
$ws.Range("D26").Value2 = "2.101.65"
$ws.Range("E26").Value2 = "  -1.29%  "
$ws.Range("D27").Formula = "'2.537"
$ws.Range("E27").Value2 = "  -5.94%  "
$ws.Range("D28").Formula = "'158.16"
$ws.Range("E28").Value2 = "  -0.45%  "
$ws.Range("D29").Formula = "'20.75"
$ws.Range("E29").Value2 = "  -2.36%  "
$ws.Range("D30").Formula = "'127.01"
$ws.Range("E30").Value2 = "  -1.51%  "
$ws.Range("D31").Formula = "'1.068"
$ws.Range("E31").Value2 = "  -3.82%  "
$ws.Range("D32").Formula = "'0.1050"
$ws.Range("E32").Value2 = "  -3.19%  "
$ws.Range("D33").Formula = "'5.573"
$ws.Range("E33").Value2 = "  -3.08%  "
$ws.Range("D34").Formula = "'3.590"
$ws.Range("E34").Value2 = "  -1.32%  "
$ws.Range("D35").Formula = "'9.388"
$ws.Range("E35").Value2 = "  -5.43%  "
$ws.Range("D36").Formula = "'0.06543"
$ws.Range("E36").Value2 = "  -3.37%  "
$ws.Range("D37").Formula = "'0.02394"
$ws.Range("E37").Value2 = "  -1.65%  "
$ws.Range("D38").Formula = "'1.316"
$ws.Range("E38").Value2 = "  +10.60%  "
$ws.Range("D39").Formula = "'0.2177"
$ws.Range("E39").Value2 = "  -2.43%  "
$ws.Range("D40").Formula = "'1.208"
$ws.Range("E40").Value2 = "  -4.70%  "
$ws.Range("D41").Formula = "'0.6393"
$ws.Range("E41").Value2 = "  -1.20%  "
$ws.Range("D42").Formula = "'11.51"
$ws.Range("E42").Value2 = "  -2.71%  "
$ws.Range("D43").Formula = "'4.927"
$ws.Range("E43").Value2 = "  -3.29%  "
$ws.Range("D44").Formula = "'1.002"
$ws.Range("E44").Value2 = "  +0.13%  "
$ws.Range("D45").Formula = "'13.37"
$ws.Range("E45").Value2 = "  -1.68%  "
$ws.Range("D46").Formula = "'0.6014"
$ws.Range("E46").Value2 = "  -1.40%  "
$ws.Range("D47").Formula = "'1.300"
$ws.Range("E47").Value2 = "  +1.07%  "
$ws.Range("D48").Formula = "'3.677"
$ws.Range("E48").Value2 = "  -2.05%  "
$ws.Range("D49").Formula = "'1.988"
$ws.Range("E49").Value2 = "  -2.15%  "
$ws.Range("D50").Formula = "'1.201"
$ws.Range("E50").Value2 = "  -0.71%  "
$ws.Range("D51").Formula = "'120.48"
$ws.Range("E51").Value2 = "  -3.92%  "
